$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new sheet after the existing one
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "parameterizedSearchTest"

# Populate header row
$ws2.Range("A1").Value = "region"
$ws2.Range("B1").Value = "country"
$ws2.Range("C1").Value = "initialPort"
$ws2.Range("D1").Value = "minCruiseLength"
$ws2.Range("E1").Value = "date"
$ws2.Range("F1").Value = "priceLimit"

# Rows 2-4, columns A-C first (C2/C3 temporarily swapped to match authoring order)
$ws2.Range("A2").Value = "BritishIslands"
$ws2.Range("B2").Value = "GreatBritain"
$ws2.Range("C2").Value = "Genoa"

$ws2.Range("A3").Value = "NearEast"
$ws2.Range("B3").Value = "Israel"
$ws2.Range("C3").Value = "Marseille"

$ws2.Range("A4").Value = "Africa"
$ws2.Range("B4").Value = "AnyCountry"
$ws2.Range("C4").Value = "Rome"

# Columns D-F for rows 2-4
$ws2.Range("D2").Value = 7
$ws2.Range("E2").Value = "20 июнь"
$ws2.Range("F2").Value = 60000

$ws2.Range("D3").Value = 7
$ws2.Range("E3").Value = "20 июнь"
$ws2.Range("F3").Value = 60000

$ws2.Range("D4").Value = 7
$ws2.Range("E4").Value = "20 июнь"
$ws2.Range("F4").Value = 100000

# Final swap of C2/C3 to correct values
$ws2.Range("C2").Value = "Marseille"
$ws2.Range("C3").Value = "Genoa"

Write-Host "done"
